$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "PO Forecast"

# --- Header row ---
$hdr = New-Object 'object[,]' 1,4
$hdr[0,0] = "ds"
$hdr[0,1] = "PO_Forecast"
$hdr[0,2] = "yhat_lower"
$hdr[0,3] = "yhat_upper"
$newSheet.Range("A1:D1").Value = $hdr

# --- Data rows (57 rows, columns A-D) ---
$data = New-Object 'object[,]' 57,4
$data[0,0] = 44948.99999999999
$data[0,1] = 162
$data[0,2] = -69.14117710873968
$data[0,3] = 389.7479749808867
$data[1,0] = 44962.99999999999
$data[1,1] = 163
$data[1,2] = -50.96300895755991
$data[1,3] = 404.7121643650302
$data[2,0] = 44969.99999999999
$data[2,1] = 164
$data[2,2] = -55.28103892213954
$data[2,3] = 371.3111155756816
$data[3,0] = 44976.99999999999
$data[3,1] = 164
$data[3,2] = -65.16269358623735
$data[3,3] = 385.3293080331908
$data[4,0] = 44983.99999999999
$data[4,1] = 165
$data[4,2] = -44.19701946219981
$data[4,3] = 377.2841418952793
$data[5,0] = 44990.99999999999
$data[5,1] = 165
$data[5,2] = -68.79404222402501
$data[5,3] = 396.1867950836595
$data[6,0] = 44997.99999999999
$data[6,1] = 166
$data[6,2] = -68.38345298022746
$data[6,3] = 390.747356212882
$data[7,0] = 45004.99999999999
$data[7,1] = 166
$data[7,2] = -60.9434435109778
$data[7,3] = 398.9370844924284
$data[8,0] = 45011.99999999999
$data[8,1] = 167
$data[8,2] = -58.36630106120064
$data[8,3] = 398.7885152410263
$data[9,0] = 45018.99999999999
$data[9,1] = 167
$data[9,2] = -71.93335672738688
$data[9,3] = 387.2447989066105
$data[10,0] = 45025.99999999999
$data[10,1] = 168
$data[10,2] = -56.53692527638101
$data[10,3] = 387.8516425691544
$data[11,0] = 45032.99999999999
$data[11,1] = 169
$data[11,2] = -75.02153497499414
$data[11,3] = 387.2292171710479
$data[12,0] = 45039.99999999999
$data[12,1] = 169
$data[12,2] = -51.81791494689263
$data[12,3] = 390.1209237524375
$data[13,0] = 45046.99999999999
$data[13,1] = 170
$data[13,2] = -53.18376261221403
$data[13,3] = 395.1623823702468
$data[14,0] = 45067.99999999999
$data[14,1] = 171
$data[14,2] = -56.49876892246955
$data[14,3] = 393.359326455125
$data[15,0] = 45088.99999999999
$data[15,1] = 173
$data[15,2] = -65.0762595809101
$data[15,3] = 389.2137246853135
$data[16,0] = 45095.99999999999
$data[16,1] = 174
$data[16,2] = -46.86630200640397
$data[16,3] = 398.0582137011378
$data[17,0] = 45116.99999999999
$data[17,1] = 175
$data[17,2] = -55.3932248491411
$data[17,3] = 407.9894175335934
$data[18,0] = 45123.99999999999
$data[18,1] = 176
$data[18,2] = -76.85674578729373
$data[18,3] = 399.5215707919106
$data[19,0] = 45130.99999999999
$data[19,1] = 176
$data[19,2] = -45.79362771027505
$data[19,3] = 384.6662860241313
$data[20,0] = 45137.99999999999
$data[20,1] = 177
$data[20,2] = -53.80462914069501
$data[20,3] = 404.2706023350705
$data[21,0] = 45144.99999999999
$data[21,1] = 177
$data[21,2] = -40.0184162090435
$data[21,3] = 390.2329644175252
$data[22,0] = 45165.99999999999
$data[22,1] = 179
$data[22,2] = -45.59217472807192
$data[22,3] = 391.1668950892308
$data[23,0] = 45186.99999999999
$data[23,1] = 181
$data[23,2] = -47.34089955586433
$data[23,3] = 398.5304907333594
$data[24,0] = 45200.99999999999
$data[24,1] = 182
$data[24,2] = -52.11060072923312
$data[24,3] = 408.9944334224615
$data[25,0] = 45214.99999999999
$data[25,1] = 183
$data[25,2] = -43.79109692934537
$data[25,3] = 416.0614421917281
$data[26,0] = 45221.99999999999
$data[26,1] = 183
$data[26,2] = -53.42880982175542
$data[26,3] = 397.2875980792705
$data[27,0] = 45270.99999999999
$data[27,1] = 187
$data[27,2] = -37.05838284784337
$data[27,3] = 418.5790075566449
$data[28,0] = 45396.99999999999
$data[28,1] = 197
$data[28,2] = -37.71778975532214
$data[28,3] = 409.29140820543
$data[29,0] = 45403.99999999999
$data[29,1] = 198
$data[29,2] = -33.13806149139964
$data[29,3] = 422.7435188059849
$data[30,0] = 45410.99999999999
$data[30,1] = 198
$data[30,2] = -13.0294941564837
$data[30,3] = 426.8720524297792
$data[31,0] = 45417.99999999999
$data[31,1] = 199
$data[31,2] = -29.3013879268508
$data[31,3] = 433.9803662952361
$data[32,0] = 45424.99999999999
$data[32,1] = 199
$data[32,2] = -29.79267094274116
$data[32,3] = 440.2444923551952
$data[33,0] = 45431.99999999999
$data[33,1] = 200
$data[33,2] = -46.48695404441958
$data[33,3] = 428.2550794069162
$data[34,0] = 45438.99999999999
$data[34,1] = 200
$data[34,2] = -33.56933055952182
$data[34,3] = 450.0053379971378
$data[35,0] = 45445.99999999999
$data[35,1] = 201
$data[35,2] = -27.38905182285995
$data[35,3] = 425.8460283406527
$data[36,0] = 45452.99999999999
$data[36,1] = 201
$data[36,2] = -10.25356167449192
$data[36,3] = 435.8483769323495
$data[37,0] = 45459.99999999999
$data[37,1] = 202
$data[37,2] = -21.34787724938449
$data[37,3] = 425.6045398307986
$data[38,0] = 45466.99999999999
$data[38,1] = 203
$data[38,2] = -15.84913325767039
$data[38,3] = 427.0864867845486
$data[39,0] = 45494.99999999999
$data[39,1] = 205
$data[39,2] = -31.59311957043742
$data[39,3] = 420.271315551419
$data[40,0] = 45501.99999999999
$data[40,1] = 205
$data[40,2] = -24.7327670613058
$data[40,3] = 421.8217645025487
$data[41,0] = 45508.99999999999
$data[41,1] = 206
$data[41,2] = -12.01408114974858
$data[41,3] = 437.3020860578099
$data[42,0] = 45515.99999999999
$data[42,1] = 206
$data[42,2] = -18.00379955847925
$data[42,3] = 444.5756037591564
$data[43,0] = 45529.99999999999
$data[43,1] = 208
$data[43,2] = -4.459946445159673
$data[43,3] = 433.9780003438468
$data[44,0] = 45536.99999999999
$data[44,1] = 208
$data[44,2] = -14.69494599566723
$data[44,3] = 429.9975219634454
$data[45,0] = 45543.99999999999
$data[45,1] = 209
$data[45,2] = -18.30465306057593
$data[45,3] = 429.9647016005885
$data[46,0] = 45557.99999999999
$data[46,1] = 210
$data[46,2] = -5.053574078026948
$data[46,3] = 410.4182053570485
$data[47,0] = 45564.99999999999
$data[47,1] = 210
$data[47,2] = -22.98975438411694
$data[47,3] = 442.6013941587495
$data[48,0] = 45571.99999999999
$data[48,1] = 211
$data[48,2] = -1.588885880928661
$data[48,3] = 448.0871982041684
$data[49,0] = 45578.99999999999
$data[49,1] = 211
$data[49,2] = -17.78489954880907
$data[49,3] = 418.8897686169318
$data[50,0] = 45585.99999999999
$data[50,1] = 212
$data[50,2] = 12.97720544991339
$data[50,3] = 451.3595867638123
$data[51,0] = 45592.99999999999
$data[51,1] = 212
$data[51,2] = -11.09200942605714
$data[51,3] = 434.6418740570913
$data[52,0] = 45599.99999999999
$data[52,1] = 213
$data[52,2] = -18.79279106612192
$data[52,3] = 435.6406881157395
$data[53,0] = 45606.99999999999
$data[53,1] = 214
$data[53,2] = -11.49256630960303
$data[53,3] = 433.0532931937977
$data[54,0] = 45613.99999999999
$data[54,1] = 214
$data[54,2] = -6.26142570368485
$data[54,3] = 455.9010043134966
$data[55,0] = 45620.99999999999
$data[55,1] = 215
$data[55,2] = -21.47882801177926
$data[55,3] = 431.7768522254316
$data[56,0] = 45627.99999999999
$data[56,1] = 215
$data[56,2] = -5.061196340132772
$data[56,3] = 434.2879377482606

$newSheet.Range("A2:D58").Value = $data

# --- Copy cell formatting from the "Weekly Quantity" sheet so styles match ---
# Header style (bold, centered, bordered) -> row 1 across A:D
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Date style (custom datetime number format) -> column A, rows 2-58
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A58").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
